# Design Brief.docx - add reviewer comments (Ethan Dawkins) to various
# spans of text throughout the document, per the commit:
#   "Linker fixed, debug + release (for real this time)"
#
# Word's COM model auto-generates comments.xml / commentsExtended.xml /
# commentsExtensible.xml / commentsIds.xml / people.xml plus the required
# relationships the first time a comment is added, so we only need to
# drive Range + Comments.Add().

$d = $word.ActiveDocument

# All new comments are authored by the same reviewer.
$word.UserName = "Ethan Dawkins"
$word.UserInitials = "ED"

# Running cursor (character offset) used to walk forward through the
# document story so that repeated words (e.g. "algorithms") resolve to
# the correct, later occurrence instead of the first one in the body.
$script:cursor = 0

function Find-NextRange([string]$text) {
    $searchRange = $d.Range($script:cursor, $d.Content.End)
    $ok = $searchRange.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find text: $text"
    }
    $script:cursor = $searchRange.End
    return $d.Range($searchRange.Start, $searchRange.End)
}

function Add-ReviewComment([string]$text, [string]$commentText) {
    $target = Find-NextRange $text
    $d.Comments.Add($target, $commentText) | Out-Null
}

# --- Comment 0: "Objective" heading -----------------------------------
Add-ReviewComment "Objective" "Confident in objective!"

# --- Comment 1: "not" in "includes, but is not limited to:" ----------
Add-ReviewComment "not" "Be precise in what the objective is"

# --- Comment 2: "useful" in "...but how is this useful? Number theory" -
Add-ReviewComment "useful" "Be confident in how the objective will be achieved"

# --- Comment 3: "but the main algorithms that " -----------------------
Add-ReviewComment "but the main algorithms that " "I will be using..."

# --- Comment 4: start of the "Cosine interpolation is ..." paragraph --
# (the comment range starts at the very beginning of that paragraph,
# i.e. before the two floating pictures anchored to it, and ends right
# after "Cosine interpolation ").
$coseParaStart = $null
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.StartsWith("Cosine interpolation")) {
        $coseParaStart = $p.Range.Start
        break
    }
}
if ($null -eq $coseParaStart) {
    throw "Could not locate the 'Cosine interpolation' paragraph"
}
$script:cursor = $coseParaStart
Add-ReviewComment "Cosine interpolation " "Mathematical operation, not advanced algorithm"

# --- Comment 5: "Universal data packaging " ---------------------------
Add-ReviewComment "Universal data packaging " "Explain exactly how this 'algorithm' works"

# --- Comment 6: "Integration" heading ----------------------------------
Add-ReviewComment "Integration" "More!!"

# --- Comment 7: "system to be integrated into " ------------------------
Add-ReviewComment "system to be integrated into " "Explain in more depth how it will be integrated"

# --- Comment 8: "++ project via Linker " --------------------------------
Add-ReviewComment "++ project via Linker " "Not just you need to link it. But also how the user will use the library, i.e. helper functions."

# --- Comment 9: "Modularity" heading -------------------------------------
Add-ReviewComment "Modularity" "How will you customize it!!"

Write-Output "Added $($d.Comments.Count) comments"
